$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(12.21, 0.16, 25.24234938621521, 12.03125),
    @(12.21, 0.16, 65.02928471565247, 30.671875),
    @(10.25, 0.16, 21.04346394538879, 8.671875)
)

$startRow = 295
for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $startRow + $i
    $rowData = $data[$i]
    $ws.Cells.Item($row, 1).Value = $rowData[0]
    $ws.Cells.Item($row, 2).Value = $rowData[1]
    $ws.Cells.Item($row, 3).Value = $rowData[2]
    $ws.Cells.Item($row, 4).Value = $rowData[3]
}
